# FBSS Quarterly Financials update
# Inserts two new quarterly columns (31-Dec-18 and 30-Sep-18) before the
# existing data (old column D shifts to F, E->G, ... K->M) and fills in the
# new quarter's figures, matching "Doing Updates for Financials".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at D:E; Excel shifts existing D:K data to F:M.
$ws.Columns("D:E").Insert()

# The new D:E columns come in unformatted; copy the number/date formatting
# that column F (the old column D) carries so the new quarters look the same
# as the rest of the table (date format in the header row, thousands format
# for the data rows).
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New values for the two inserted quarter columns, keyed by row number:
# @(D-column value, E-column value). Rows not listed here were blank in the
# original table and stay blank.
$newQuarterData = @{
    7 = @(43465, 43373)
    8 = @(7100, 6700)
    9 = @("NA", "NA")
    10 = @("NA", "NA")
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    17 = @(1000, 1000)
    18 = @(6100, 5700)
    20 = @(-4400, -4200)
    21 = @(2000, 1800)
    22 = @(0, 0)
    23 = @(1700, 1500)
    24 = @(100, 200)
    25 = @(0, 0)
    26 = @(1600, 1300)
    27 = @(1600, 1300)
    28 = @(0, 0)
    29 = @(0, "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(4400, 4200)
    33 = @(1600, 1300)
    34 = @(0, 0)
    35 = @(1600, 1300)
    38 = @(43465, 43373)
    41 = @(6800, 5400)
    42 = @(62500, 18100)
    43 = @(0, 0)
    44 = @(0, 0)
    45 = @(0, 0)
    46 = @(0, 0)
    47 = @(0, 0)
    48 = @(18200, 18300)
    49 = @(0, 0)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(0, 0)
    53 = @(0, 0)
    54 = @(730800, 670000)
    57 = @(0, 0)
    58 = @(0, 0)
    59 = @(0, 0)
    60 = @(0, 0)
    61 = @(4100, 4100)
    62 = @(0, 0)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(670800, 611700)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(44800, 43700)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(60000, 58300)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(1600, 1300)
    83 = @(300, 300)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(2900, 2900)
    91 = @(-200, -100)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-14700, -23600)
    96 = @(-500, -500)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(57900, 16900)
    101 = @(0, 0)
    102 = @(46100, -3700)
}

foreach ($r in $newQuarterData.Keys) {
    $pair = $newQuarterData[$r]
    if ($null -ne $pair[0]) {
        $ws.Cells.Item($r, 4).Value = $pair[0]
    }
    if ($null -ne $pair[1]) {
        $ws.Cells.Item($r, 5).Value = $pair[1]
    }
}
